# Auto-generated Excel COM-interop script to apply scheduled runner price/profit updates
# across sheets ALC, ARM, CRP, CUL, GSM, LTW, WVR per diff of Omega_Profits workbook.
$wb = $excel.ActiveWorkbook

# Sheet ALC, Row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1020.62164
$ws.Range("I15").Value = 1020.62164
$ws.Range("K15").Value = 3061.86492
$ws.Range("M15").Value = -2892.86492

# Sheet ALC, Row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 10597.906
$ws.Range("I76").Value = 10696.174
$ws.Range("J76").Value = 10346.777
$ws.Range("K76").Value = 10696.174
$ws.Range("L76").Value = 10346.777
$ws.Range("M76").Value = -10381.174
$ws.Range("N76").Value = -10976.777

# Sheet ALC, Row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 10597.906
$ws.Range("I79").Value = 10696.174
$ws.Range("J79").Value = 10346.777
$ws.Range("K79").Value = 10696.174
$ws.Range("L79").Value = 10346.777
$ws.Range("M79").Value = -9604.174000000001
$ws.Range("N79").Value = -12530.777

# Sheet ALC, Row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1678.7241
$ws.Range("J137").Value = 2063.2222
$ws.Range("L137").Value = 6189.6666
$ws.Range("N137").Value = -11289.6666

# Sheet ALC, Row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2959.415
$ws.Range("I138").Value = 1555.55
$ws.Range("J138").Value = 3810.2424
$ws.Range("K138").Value = 4666.65
$ws.Range("L138").Value = 11430.7272
$ws.Range("M138").Value = 473.3500000000004
$ws.Range("N138").Value = -21710.7272

# Sheet ALC, Row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 75000
$ws.Range("J139").Value = 75000
$ws.Range("L139").Value = 75000
$ws.Range("N139").Value = -85280

# Sheet ARM, Row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5354.184
$ws.Range("I132").Value = 5321
$ws.Range("K132").Value = 15963
$ws.Range("M132").Value = -13433

# Sheet ARM, Row 137
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 79999.5
$ws.Range("J137").Value = 79999
$ws.Range("L137").Value = 79999
$ws.Range("N137").Value = -90199

# Sheet ARM, Row 138
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 88391
$ws.Range("J138").Value = 88391
$ws.Range("L138").Value = 88391
$ws.Range("N138").Value = -98671

# Sheet ARM, Row 141
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 89999.75
$ws.Range("J141").Value = 89999.75
$ws.Range("L141").Value = 89999.75
$ws.Range("N141").Value = -100359.75

# Sheet CRP, Row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 809.8889
$ws.Range("I22").Value = 548.3333
$ws.Range("K22").Value = 548.3333
$ws.Range("M22").Value = -198.3333

# Sheet CRP, Row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4790.1113
$ws.Range("I58").Value = 4992.3184
$ws.Range("K58").Value = 4992.3184
$ws.Range("M58").Value = -4789.3184

# Sheet CRP, Row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6248.773
$ws.Range("I132").Value = 5120.815
$ws.Range("J132").Value = 8040.2354
$ws.Range("K132").Value = 15362.445
$ws.Range("L132").Value = 24120.7062
$ws.Range("M132").Value = -12832.445
$ws.Range("N132").Value = -29180.7062

# Sheet CRP, Row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 7036.875
$ws.Range("I134").Value = 6402.096
$ws.Range("K134").Value = 19206.288
$ws.Range("M134").Value = -16671.288

# Sheet CRP, Row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4790.1113
$ws.Range("I136").Value = 4992.3184
$ws.Range("K136").Value = 14976.9552
$ws.Range("M136").Value = -12426.9552

# Sheet CUL, Row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 640.5
$ws.Range("I5").Value = 640.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1921.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1809.5
$ws.Range("N5").ClearContents()

# Sheet CUL, Row 32
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 8000
$ws.Range("I32").Value = 8000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 24000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -23717
$ws.Range("N32").ClearContents()

# Sheet CUL, Row 61
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 528.4286
$ws.Range("J61").Value = 499.5
$ws.Range("L61").Value = 1498.5
$ws.Range("N61").Value = -1928.5

# Sheet CUL, Row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I75").Value = 3283
$ws.Range("J75").Value = 3497.8
$ws.Range("K75").Value = 9849
$ws.Range("L75").Value = 10493.4
$ws.Range("M75").Value = -8851
$ws.Range("N75").Value = -12489.4

# Sheet CUL, Row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I78").Value = 3283
$ws.Range("J78").Value = 3497.8
$ws.Range("K78").Value = 29547
$ws.Range("L78").Value = 31480.2
$ws.Range("M78").Value = -24555
$ws.Range("N78").Value = -41464.2

# Sheet CUL, Row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 974.8333
$ws.Range("I107").Value = 861.875
$ws.Range("J107").Value = 1200.75
$ws.Range("K107").Value = 2585.625
$ws.Range("L107").Value = 3602.25
$ws.Range("M107").Value = -665.625
$ws.Range("N107").Value = -7442.25

# Sheet CUL, Row 114
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2131.1667
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()

# Sheet CUL, Row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1923.1666
$ws.Range("I131").Value = 1375.5714
$ws.Range("K131").Value = 4126.7142
$ws.Range("M131").Value = 913.2857999999997

# Sheet CUL, Row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 640.5
$ws.Range("I135").Value = 640.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 5764.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3229.5
$ws.Range("N135").ClearContents()

# Sheet GSM, Row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# Sheet GSM, Row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5573.263
$ws.Range("I126").Value = 4606.364
$ws.Range("J126").Value = 6902.75
$ws.Range("K126").Value = 13819.092
$ws.Range("L126").Value = 20708.25
$ws.Range("M126").Value = -11349.092
$ws.Range("N126").Value = -25648.25

# Sheet GSM, Row 130
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 56600

# Sheet GSM, Row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6892.5
$ws.Range("I132").Value = 5848
$ws.Range("K132").Value = 17544
$ws.Range("M132").Value = -15014

# Sheet GSM, Row 135
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 84827
$ws.Range("J135").Value = 84827
$ws.Range("L135").Value = 84827
$ws.Range("N135").Value = -94967

# Sheet GSM, Row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 88333
$ws.Range("J138").Value = 88333
$ws.Range("L138").Value = 88333
$ws.Range("N138").Value = -98613

# Sheet LTW, Row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 869.9231
$ws.Range("I16").Value = 952.9091
$ws.Range("K16").Value = 952.9091
$ws.Range("M16").Value = -782.9091

# Sheet LTW, Row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5000
$ws.Range("I22").Value = 5000
$ws.Range("K22").Value = 5000
$ws.Range("M22").Value = -4705

# Sheet LTW, Row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 5000
$ws.Range("I27").Value = 5000
$ws.Range("K27").Value = 5000
$ws.Range("M27").Value = -4893

# Sheet LTW, Row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 10777.25
$ws.Range("I122").Value = 10777.25
$ws.Range("K122").Value = 32331.75
$ws.Range("M122").Value = -29881.75

# Sheet WVR, Row 41
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11999.5
$ws.Range("I41").Value = 11999
$ws.Range("K41").Value = 11999
$ws.Range("M41").Value = -11609

# Sheet WVR, Row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 961.5454999999999
$ws.Range("I113").Value = 961.5454999999999
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2884.6365
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -714.6364999999996
$ws.Range("N113").ClearContents()

# Sheet WVR, Row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3696.5854
$ws.Range("I136").Value = 3699.4375
$ws.Range("K136").Value = 11098.3125
$ws.Range("M136").Value = -8548.3125

# Sheet WVR, Row 138
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280
